# "end of sept updates" - update flowchart box counts/labels and
# reposition two of the boxes on the single slide of the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Rectangle 4 (shape id 5): "Women with complete biological age
# biomarker panel" -> "Complete biological age biomarker panel",
# n = 8,006 -> n = 5,870, and move it down to where Rectangle 5 used
# to sit.
$rect4 = $s.Shapes.Item("Rectangle 4")
$rect4.Top = 2275818 / 12700
$tr4 = $rect4.TextFrame.TextRange
$tr4.Runs(1, 1).Text = "Complete biological age biomarker panel"
$tr4.Runs(3, 1).Text = " = 5,870"

# --- Rectangle 5 (shape id 6): "Ages 18-84 and not currently
# pregnant" -> "Women ages 18-84 and not currently pregnant",
# n = 3,651 -> n = 13,929, and move it up to where Rectangle 4 used
# to sit.
$rect5 = $s.Shapes.Item("Rectangle 5")
$rect5.Left = 2244566 / 12700
$rect5.Top = 1420496 / 12700
$tr5 = $rect5.TextFrame.TextRange
$tr5.Runs(1, 1).Text = "Women ages 18-84 and not currently pregnant"
$tr5.Runs(3, 1).Text = " = 13,929"

# --- Rectangle 6 (shape id 7): "Complete covariate information",
# n = 2,696 -> n = 4,418.
$rect6 = $s.Shapes.Item("Rectangle 6")
$tr6 = $rect6.TextFrame.TextRange
$tr6.Runs(3, 1).Text = " = 4,418"

# --- Rectangle 9 (shape id 10): "Zero to six live births reported"
# -> "Zero to seven live births reported", n = 3,235 -> n = 5,184.
$rect9 = $s.Shapes.Item("Rectangle 9")
$tr9 = $rect9.TextFrame.TextRange
$tr9.Runs(1, 1).Text = "Zero to seven live births reported"
$tr9.Runs(3, 1).Text = " = 5,184"

# --- Rectangle 14 (shape id 15): "Data on years since last birth",
# n = 2,056 -> n = 3,587.
$rect14 = $s.Shapes.Item("Rectangle 14")
$tr14 = $rect14.TextFrame.TextRange
$tr14.Runs(3, 1).Text = " = 3,587"
